$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.498.61'
$ws.Range("E2").Value = '  +2.52%  '
$ws.Range("D3").Value = '2.191.80'
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '253.54'
$ws.Range("E5").Value = '  +6.36%  '
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("D7").Value = '73.66'
$ws.Range("E7").Value = '  +2.83%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.588'
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("D10").Value = '40.04'
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("D11").Value = '0.0916'
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("D12").Value = '6.81'
$ws.Range("E12").Value = '  +2.32%  '
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").Value = '2.523.71'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '14.34'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '2.195.09'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("D17").Value = '0.775'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '42.411.87'
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").Value = '71.06'
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("D21").Value = '5.91'
$ws.Range("E21").Value = '  +2.73%  '
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").Value = '2.20'
$ws.Range("E22").Value = '  +9.03%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '9.69'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").Value = '227.92'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '3.35'
$ws.Range("E27").Value = '  +1.57%  '
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("D29").Value = '2.14'
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").Value = '37.05'
$ws.Range("E30").Value = '  +12.15%  '
$ws.Range("D31").Value = '168.89'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '20.05'
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").Value = '0.0804'
$ws.Range("E33").Value = '  +5.04%  '
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").Value = '0.121'
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("D37").Value = '4.35'
$ws.Range("E37").Value = '  +1.72%  '
$ws.Range("D38").Value = '0.0332'
$ws.Range("E38").Value = '  +10.07%  '
$ws.Range("D39").Value = '12.10'
$ws.Range("E39").Value = '  +2.02%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("E41").Value = '  +5.24%  '
$ws.Range("D42").Value = '5.21'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").Value = '59.03'
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("D44").Value = '102.62'
$ws.Range("E44").Value = '  +7.15%  '
$ws.Range("D45").Value = '0.469'
$ws.Range("E45").Value = '  +19.20%  '
$ws.Range("D46").Value = '8.32'
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").Value = '0.0974'
$ws.Range("E47").Value = '  +1.81%  '
$ws.Range("D48").Value = '2.41'
$ws.Range("E48").Value = '  +11.11%  '
$ws.Range("E49").Value = '  +2.78%  '
$ws.Range("D50").Value = '1.13'
$ws.Range("E50").Value = '  +1.96%  '
